# prob28 in p&s part 2: expand the "치역(range)" building-block breakdown
# from 3 rows (x0016-x0018) into 6 rows (x0016-x0021), adding a per-case
# breakdown of the possible ranges {1,2,3}/{1,2,4}/{1,3,4}/{2,3,4} tagged
# with the problem id 32111_x28 in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make room: rows below 168 (the old gap + the "y0001"/"z0001" markers at
# rows 181/191) need to end up 11 rows further down (192/202), while rows
# 166-168 themselves grow from 3 used rows to 6 (166-171). Inserting 11
# blank rows starting at 169 achieves both in one shot.
$ws.Rows("169:179").Insert()

# Row 166 - unchanged key, new (more specific) building-block text
$ws.Range("A166").Value = "x0016"
$ws.Range("B166").Value = "함수의 조건에 맞는 가능한 치역을 모두 구합니다."

# Row 167 - new key/text, tagged with the problem id
$ws.Range("A167").Value = "x0017"
$ws.Range("B167").Value = "치역이 `$\{1, 2, 3\}`$ 인 경우 조건에 맞는 함수의 개수를 구합니다."
$ws.Range("C167").Value = "32111_x28"

# Row 168
$ws.Range("A168").Value = "x0018"
$ws.Range("B168").Value = "치역이 `$\{1, 2, 4\}`$ 인 경우 조건에 맞는 함수의 개수를 구합니다."
$ws.Range("C168").Value = "32111_x28"

# Row 169 (new row)
$ws.Range("A169").Value = "x0019"
$ws.Range("B169").Value = "치역이 `$\{1, 3, 4\}`$ 인 경우 조건에 맞는 함수의 개수를 구합니다."
$ws.Range("C169").Value = "32111_x28"

# Row 170 (new row)
$ws.Range("A170").Value = "x0020"
$ws.Range("B170").Value = "치역이 `$\{2, 3, 4\}`$ 인 경우 조건에 맞는 함수의 개수를 구합니다."
$ws.Range("C170").Value = "32111_x28"

# Row 171 (new row) - carries forward the old final summary text (x0018 -> x0021)
$ws.Range("A171").Value = "x0021"
$ws.Range("B171").Value = "각각의 개수를 모두 더해서 조건을 만족시키는 전체 개수를 구합니다. "

# Update the sheet view to match the edited location
$ws.Application.ActiveWindow.ScrollRow = 157
$ws.Range("B179").Select()
